$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H11").Value = 95.14286
$ws.Range("I11").Value = 95.14286
$ws.Range("K11").Value = 95.14286
$ws.Range("M11").Value = 44.85714
$ws.Range("H80").Value = 4803348.5
$ws.Range("I80").Value = 3269097.8
$ws.Range("K80").Value = 9807293.399999999
$ws.Range("M80").Value = -9806295.399999999
$ws.Range("H83").Value = 4803348.5
$ws.Range("I83").Value = 3269097.8
$ws.Range("K83").Value = 29421880.2
$ws.Range("M83").Value = -29416888.2
$ws.Range("H97").Value = 4540.6665
$ws.Range("J97").Value = 4540.6665
$ws.Range("L97").Value = 13621.9995
$ws.Range("N97").Value = -14613.9995
$ws.Range("H98").Value = 1538.5834
$ws.Range("I98").Value = 1416.4
$ws.Range("K98").Value = 1416.4
$ws.Range("M98").Value = 81.59999999999991
$ws.Range("H122").Value = 1538.5834
$ws.Range("I122").Value = 1416.4
$ws.Range("K122").Value = 4249.200000000001
$ws.Range("M122").Value = -1799.200000000001
$ws.Range("H131").Value = 2530002.5
$ws.Range("I131").Value = 812
$ws.Range("K131").Value = 2436
$ws.Range("M131").Value = 2604
$ws.Range("H132").Value = 3429.8704
$ws.Range("I132").Value = 2009.6171
$ws.Range("K132").Value = 6028.8513
$ws.Range("M132").Value = -3498.8513
$ws.Range("H135").Value = 2477.42
$ws.Range("I135").Value = 1835.8387
$ws.Range("K135").Value = 16522.5483
$ws.Range("M135").Value = -13987.5483
$ws.Range("H137").Value = 1620.9032
$ws.Range("I137").Value = 1233.96
$ws.Range("J137").Value = 3233.1667
$ws.Range("K137").Value = 3701.88
$ws.Range("L137").Value = 9699.500100000001
$ws.Range("M137").Value = -1151.88
$ws.Range("N137").Value = -14799.5001
$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("H30").Value = 27096.334
$ws.Range("I30").Value = 650
$ws.Range("K30").Value = 650
$ws.Range("M30").Value = -500
$ws.Range("H46").Value = 31710
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 31710
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 31710
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -32348
$ws.Range("H97").Value = 1483
$ws.Range("I97").Value = 742.5789
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 742.5789
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = -246.5789
$ws.Range("N97").Value = -5992
$ws.Range("H102").Value = 26317854
$ws.Range("I102").Value = 31251638
$ws.Range("K102").Value = 31251638
$ws.Range("M102").Value = -31250016
$ws.Range("H133").Value = 67249.75
$ws.Range("J133").Value = 67249.75
$ws.Range("L133").Value = 67249.75
$ws.Range("N133").Value = -72309.75
$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("H105").Value = 718324.75
$ws.Range("I105").Value = 1272677
$ws.Range("J105").Value = 5586.143
$ws.Range("K105").Value = 1272677
$ws.Range("L105").Value = 5586.143
$ws.Range("M105").Value = -1270930
$ws.Range("N105").Value = -9080.143
$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H7").Value = 598.125
$ws.Range("I7").Value = 446.75
$ws.Range("J7").Value = 749.5
$ws.Range("K7").Value = 446.75
$ws.Range("L7").Value = 749.5
$ws.Range("M7").Value = -333.75
$ws.Range("N7").Value = -975.5
$ws.Range("H22").Value = 702.7143
$ws.Range("I22").Value = 605.375
$ws.Range("K22").Value = 605.375
$ws.Range("M22").Value = -255.375
$ws.Range("H31").Value = 2888.5833
$ws.Range("I31").Value = 2889.8518
$ws.Range("J31").Value = 2884.7778
$ws.Range("K31").Value = 2889.8518
$ws.Range("L31").Value = 2884.7778
$ws.Range("M31").Value = -2594.8518
$ws.Range("N31").Value = -3474.7778
$ws.Range("H34").Value = 2888.5833
$ws.Range("I34").Value = 2889.8518
$ws.Range("J34").Value = 2884.7778
$ws.Range("K34").Value = 2889.8518
$ws.Range("L34").Value = 2884.7778
$ws.Range("M34").Value = -2687.8518
$ws.Range("N34").Value = -3288.7778
$ws.Range("H99").Value = 8210.916999999999
$ws.Range("I99").Value = 4553.1
$ws.Range("J99").Value = 26500
$ws.Range("K99").Value = 4553.1
$ws.Range("L99").Value = 26500
$ws.Range("M99").Value = -3055.1
$ws.Range("N99").Value = -29496
$ws.Range("H108").Value = 99999
$ws.Range("J108").Value = 99999
$ws.Range("L108").Value = 99999
$ws.Range("N108").Value = -107679
$ws.Range("H126").Value = 8210.916999999999
$ws.Range("I126").Value = 4553.1
$ws.Range("J126").Value = 26500
$ws.Range("K126").Value = 13659.3
$ws.Range("L126").Value = 79500
$ws.Range("M126").Value = -11189.3
$ws.Range("N126").Value = -84440
$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H107").Value = 2531977.8
$ws.Range("J107").Value = 3796510.5
$ws.Range("L107").Value = 11389531.5
$ws.Range("N107").Value = -11393371.5
$ws.Range("H118").Value = 7812.4
$ws.Range("I118").Value = 4976.778
$ws.Range("K118").Value = 14930.334
$ws.Range("M118").Value = -13687.334
$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H57").Value = 37500
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 37500
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 37500
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -39140
$ws.Range("H64").Value = 60000.555
$ws.Range("J64").Value = 60000.555
$ws.Range("L64").Value = 60000.555
$ws.Range("N64").Value = -60496.555
$ws.Range("H67").Value = 60000.555
$ws.Range("J67").Value = 60000.555
$ws.Range("L67").Value = 60000.555
$ws.Range("N67").Value = -61716.555
$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H46").Value = 1352.1111
$ws.Range("I46").Value = 1309.8572
$ws.Range("K46").Value = 1309.8572
$ws.Range("M46").Value = -1121.8572
$ws.Range("H55").Value = 1154.45
$ws.Range("I55").Value = 851.63635
$ws.Range("K55").Value = 851.63635
$ws.Range("M55").Value = -678.63635
$ws.Range("H93").Value = 1987379
$ws.Range("I93").Value = 1491.5625
$ws.Range("K93").Value = 1491.5625
$ws.Range("M93").Value = -243.5625
$ws.Range("H132").Value = 3841.0476
$ws.Range("I132").Value = 3155.6
$ws.Range("J132").Value = 5554.6665
$ws.Range("K132").Value = 9466.799999999999
$ws.Range("L132").Value = 16663.9995
$ws.Range("M132").Value = -6936.799999999999
$ws.Range("N132").Value = -21723.9995
$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H82").Value = 99999
$ws.Range("J82").Value = 99999
$ws.Range("L82").Value = 99999
$ws.Range("N82").Value = -100765
$ws.Range("H85").Value = 99999
$ws.Range("J85").Value = 99999
$ws.Range("L85").Value = 99999
$ws.Range("N85").Value = -102651
$ws.Range("H100").Value = 579.2353000000001
$ws.Range("I100").Value = 575
$ws.Range("K100").Value = 1150
$ws.Range("M100").Value = -609
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("H122").Value = 2509.4614
$ws.Range("I122").Value = 2293.4546
$ws.Range("K122").Value = 6880.3638
$ws.Range("M122").Value = -4430.3638
$ws.Range("H132").Value = 3407.5
$ws.Range("I132").Value = 3212.4666
$ws.Range("J132").Value = 3992.6
$ws.Range("K132").Value = 9637.399800000001
$ws.Range("L132").Value = 11977.8
$ws.Range("M132").Value = -7107.399800000001
$ws.Range("N132").Value = -17037.8

Write-Host "Applied all changes"